$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Piloted adoption of MLFlow*") {
        $p.Range.Delete()
        break
    }
}
